$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure values are stored as text (preserve leading zeros / avoid date or
# numeric auto-conversion), matching the original inlineStr string cells.

# Row 2 updates
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "094439854"

$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = "ΤΡΑΚΑΔΑΣ Α.Ε."

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = ""

$ws.Range("H2").NumberFormat = "@"
$ws.Range("H2").Value = ""

# Row 3 updates
$ws.Range("A3").NumberFormat = "@"
$ws.Range("A3").Value = "400008195607600"

$ws.Range("F3").NumberFormat = "@"
$ws.Range("F3").Value = "2025-01-04"

$ws.Range("I3").NumberFormat = "@"
$ws.Range("I3").Value = "34.34"

$ws.Range("J3").NumberFormat = "@"
$ws.Range("J3").Value = "8.24"

$ws.Range("K3").NumberFormat = "@"
$ws.Range("K3").Value = "42.58"
